# Slots.xlsx update — swap a handful of ingredient names on the single
# "Tabelle1" sheet and move the active selection.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 holds the ingredient assigned to each "Slot N" header in row 1.
$ws.Range("B2").Value = "Gin"
$ws.Range("D2").Value = "Tonic Water"
$ws.Range("E2").Value = "Ginger Ale"
$ws.Range("F2").Value = "Blue Curacao"
$ws.Range("G2").Value = "Zitronensaft"
$ws.Range("J2").Value = "Sahne"

# Move the active cell / selection as recorded in the saved view state.
$ws.Range("L11").Select()
